# Add a new craftsman row (row 7) to the worksheet, mirroring the
# existing data pattern already present in row 6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "احمد ابو رسلان"
$ws.Range("B7").Value = "ابواب وشبابيك بلاستك"
$ws.Range("C7").Value = "شارع المعامل"
$ws.Range("D7").Value = 75000000

$ws.Range("D7").Select()
